# Commit: "Wed, Apr 01, 2020 10:05:25 AM"
#
# The authored OOXML diff swaps the contents of ppt/theme/theme1.xml
# (was the "Integral" colour theme used by the slide master / main deck)
# and ppt/theme/theme2.xml (was the "Office Theme" colour theme used by
# the notes master) - i.e. the deck's visible theme becomes the stock
# "Office" palette while the notes master keeps the old "Integral"
# palette. The <a:fontScheme> and <a:fmtScheme> blocks are byte-for-byte
# identical between the two themes already, so the only observable
# difference is the 12 theme colours (and the name="" labels, which
# PowerPoint's automation surface does not expose for editing).
#
# PowerPoint's object model only exposes the *colour* part of a theme
# (ThemeColorScheme), reachable from the slide master's Theme. Apply the
# destination ("Office") palette there, in the standard
# MsoThemeColorSchemeIndex order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink

$p = $ppt.ActivePresentation
$scheme = $p.SlideMaster.Theme.ThemeColorScheme

function HexToComRGB($r, $g, $b) {
    return $b * 65536 + $g * 256 + $r
}

# Office theme palette (target state of ppt/theme/theme1.xml)
$officeColors = @(
    (HexToComRGB 0x00 0x00 0x00),  # 1  dk1      000000
    (HexToComRGB 0xFF 0xFF 0xFF),  # 2  lt1      FFFFFF
    (HexToComRGB 0x44 0x54 0x6A),  # 3  dk2      44546A
    (HexToComRGB 0xE7 0xE6 0xE6),  # 4  lt2      E7E6E6
    (HexToComRGB 0x5B 0x9B 0xD5),  # 5  accent1  5B9BD5
    (HexToComRGB 0xED 0x7D 0x31),  # 6  accent2  ED7D31
    (HexToComRGB 0xA5 0xA5 0xA5),  # 7  accent3  A5A5A5
    (HexToComRGB 0xFF 0xC0 0x00),  # 8  accent4  FFC000
    (HexToComRGB 0x44 0x72 0xC4),  # 9  accent5  4472C4
    (HexToComRGB 0x70 0xAD 0x47),  # 10 accent6  70AD47
    (HexToComRGB 0x05 0x63 0xC1),  # 11 hlink    0563C1
    (HexToComRGB 0x95 0x4F 0x72)   # 12 folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $scheme.Item($i).RGB = $officeColors[$i - 1]
}
